$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Finished" column (H) with "X" for these rows
$rows = @(29, 31, 33, 35, 37, 39, 41, 43)
foreach ($r in $rows) {
    $ws.Range("H$r").Value = "X"
}

# Update the active selection to reflect where the user left off
$ws.Range("H44").Select()
